$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 10.95798036489552
$ws.Cells.Item(2, 3).Value = 6.778798356056
$ws.Cells.Item(2, 4).Value = 15.04880060961699
$ws.Cells.Item(2, 5).Value = 16.4729647258361
$ws.Cells.Item(2, 7).Value = 3.687444081314091
$ws.Cells.Item(2, 9).Value = 25.8528698636452
$ws.Cells.Item(2, 10).Value = 9.41152566838297
$ws.Cells.Item(2, 11).Value = 11.54045547095927
$ws.Cells.Item(2, 15).Value = 29.04172500636703
$ws.Cells.Item(3, 2).Value = 10.70636533411506
$ws.Cells.Item(3, 3).Value = 6.609163136605951
$ws.Cells.Item(3, 4).Value = 14.9888626809102
$ws.Cells.Item(3, 5).Value = 16.41335940663106
$ws.Cells.Item(3, 7).Value = 3.689656741107035
$ws.Cells.Item(3, 9).Value = 25.94345657797771
$ws.Cells.Item(3, 10).Value = 9.419584828567949
$ws.Cells.Item(3, 11).Value = 11.37637839648652
$ws.Cells.Item(3, 15).Value = 29.11715507499949
$ws.Cells.Item(4, 2).Value = 10.5509039793706
$ws.Cells.Item(4, 3).Value = 6.504108955576851
$ws.Cells.Item(4, 4).Value = 14.95535291444177
$ws.Cells.Item(4, 5).Value = 16.38039660073805
$ws.Cells.Item(4, 7).Value = 3.691087016290683
$ws.Cells.Item(4, 9).Value = 26.00386094943715
$ws.Cells.Item(4, 10).Value = 9.425947144675961
$ws.Cells.Item(4, 11).Value = 11.27652650288468
$ws.Cells.Item(4, 15).Value = 29.16939919066404
$ws.Cells.Item(5, 2).Value = 10.48740313412078
$ws.Cells.Item(5, 3).Value = 6.461138320374546
$ws.Cells.Item(5, 4).Value = 14.94253554590469
$ws.Cells.Item(5, 5).Value = 16.36788781076492
$ws.Cells.Item(5, 7).Value = 3.691687951790396
$ws.Cells.Item(5, 9).Value = 26.02967731591003
$ws.Cells.Item(5, 10).Value = 9.428895572740412
$ws.Cells.Item(5, 11).Value = 11.23610883645784
$ws.Cells.Item(5, 15).Value = 29.1921767601515
$ws.Cells.Item(6, 2).Value = 10.47685268960708
$ws.Cells.Item(6, 3).Value = 6.453995371599497
$ws.Cells.Item(6, 4).Value = 14.94045813126051
$ws.Cells.Item(6, 5).Value = 16.36586679517717
$ws.Cells.Item(6, 7).Value = 3.691788830800717
$ws.Cells.Item(6, 9).Value = 26.03403659365874
$ws.Cells.Item(6, 10).Value = 9.429406646137316
$ws.Cells.Item(6, 11).Value = 11.22941540184484
$ws.Cells.Item(6, 15).Value = 29.19604870634345
$ws.Cells.Item(7, 2).Value = 10.55004806107288
$ws.Cells.Item(7, 3).Value = 6.503529998638655
$ws.Cells.Item(7, 4).Value = 14.95517664861507
$ws.Cells.Item(7, 5).Value = 16.38022415025452
$ws.Cells.Item(7, 7).Value = 3.691095047406124
$ws.Cells.Item(7, 9).Value = 26.00420425724704
$ws.Cells.Item(7, 10).Value = 9.425985467678668
$ws.Cells.Item(7, 11).Value = 11.275980247708
$ws.Cells.Item(7, 15).Value = 29.16970035811271
$ws.Cells.Item(8, 2).Value = 10.87148018020191
$ws.Cells.Item(8, 3).Value = 6.720532776289271
$ws.Cells.Item(8, 4).Value = 15.02745714217443
$ws.Cells.Item(8, 5).Value = 16.4516645241537
$ws.Cells.Item(8, 7).Value = 3.68819216203753
$ws.Cells.Item(8, 9).Value = 25.88311020318371
$ws.Cells.Item(8, 10).Value = 9.414011128866285
$ws.Cells.Item(8, 11).Value = 11.48372565518617
$ws.Cells.Item(8, 15).Value = 29.06650057314192
$ws.Cells.Item(9, 2).Value = 11.49015743090058
$ws.Cells.Item(9, 3).Value = 7.136233297770364
$ws.Cells.Item(9, 4).Value = 15.19482502058606
$ws.Cells.Item(9, 5).Value = 16.62012416207541
$ws.Cells.Item(9, 7).Value = 3.68306577535416
$ws.Cells.Item(9, 9).Value = 25.6836831462692
$ws.Cells.Item(9, 10).Value = 9.401738602783984
$ws.Cells.Item(9, 11).Value = 11.89606372772888
$ws.Cells.Item(9, 15).Value = 28.91132257459875
$ws.Cells.Item(10, 2).Value = 11.93259591603981
$ws.Cells.Item(10, 3).Value = 7.432251781077821
$ws.Cells.Item(10, 4).Value = 15.33268066542487
$ws.Cells.Item(10, 5).Value = 16.76047493374503
$ws.Cells.Item(10, 7).Value = 3.679640799966627
$ws.Cells.Item(10, 9).Value = 25.56046329855154
$ws.Cells.Item(10, 10).Value = 9.399538784449399
$ws.Cells.Item(10, 11).Value = 12.19920254372831
$ws.Cells.Item(10, 15).Value = 28.82627665522581
$ws.Cells.Item(11, 2).Value = 12.13028393844628
$ws.Cells.Item(11, 3).Value = 7.564234943062457
$ws.Cells.Item(11, 4).Value = 15.39845280870893
$ws.Cells.Item(11, 5).Value = 16.82775330478468
$ws.Cells.Item(11, 7).Value = 3.678156015753553
$ws.Cells.Item(11, 9).Value = 25.50949013049784
$ws.Cells.Item(11, 10).Value = 9.400013068760385
$ws.Cells.Item(11, 11).Value = 12.33656749211568
$ws.Cells.Item(11, 15).Value = 28.79391298961574
$ws.Cells.Item(12, 2).Value = 12.20455563639313
$ws.Cells.Item(12, 3).Value = 7.613780359520704
$ws.Cells.Item(12, 4).Value = 15.42378307584074
$ws.Cells.Item(12, 5).Value = 16.85370747827238
$ws.Cells.Item(12, 7).Value = 3.677604239609629
$ws.Cells.Item(12, 9).Value = 25.49092042421379
$ws.Cells.Item(12, 10).Value = 9.400404158963742
$ws.Cells.Item(12, 11).Value = 12.38846014689001
$ws.Cells.Item(12, 15).Value = 28.78256974990783
$ws.Cells.Item(13, 2).Value = 12.18858718566719
$ws.Cells.Item(13, 3).Value = 7.603129892485256
$ws.Cells.Item(13, 4).Value = 15.41830917029063
$ws.Cells.Item(13, 5).Value = 16.8480968209329
$ws.Cells.Item(13, 7).Value = 3.677722609262996
$ws.Cells.Item(13, 9).Value = 25.49488711961091
$ws.Cells.Item(13, 10).Value = 9.40031053539736
$ws.Cells.Item(13, 11).Value = 12.37729041228503
$ws.Cells.Item(13, 15).Value = 28.78497210867652
$ws.Cells.Item(14, 2).Value = 12.13640651246102
$ws.Cells.Item(14, 3).Value = 7.568320032284132
$ws.Cells.Item(14, 4).Value = 15.40052834819898
$ws.Cells.Item(14, 5).Value = 16.82987909999256
$ws.Cells.Item(14, 7).Value = 3.678110411085608
$ws.Cells.Item(14, 9).Value = 25.50794769210209
$ws.Cells.Item(14, 10).Value = 9.40004100926889
$ws.Cells.Item(14, 11).Value = 12.34083949372885
$ws.Cells.Item(14, 15).Value = 28.79296147806065
$ws.Cells.Item(15, 2).Value = 12.10436561236422
$ws.Cells.Item(15, 3).Value = 7.546940122519592
$ws.Cells.Item(15, 4).Value = 15.38969177173466
$ws.Cells.Item(15, 5).Value = 16.81878188485378
$ws.Cells.Item(15, 7).Value = 3.678349314018951
$ws.Cells.Item(15, 9).Value = 25.51604315756003
$ws.Cells.Item(15, 10).Value = 9.39990343908822
$ws.Cells.Item(15, 11).Value = 12.31849465764176
$ws.Cells.Item(15, 15).Value = 28.79797406803172
$ws.Cells.Item(16, 2).Value = 11.9195980826387
$ws.Cells.Item(16, 3).Value = 7.423568239393831
$ws.Cells.Item(16, 4).Value = 15.32844251249985
$ws.Cells.Item(16, 5).Value = 16.75614582951967
$ws.Cells.Item(16, 7).Value = 3.679739303505493
$ws.Cells.Item(16, 9).Value = 25.56389692972532
$ws.Cells.Item(16, 10).Value = 9.399537419155218
$ws.Cells.Item(16, 11).Value = 12.1902105611477
$ws.Cells.Item(16, 15).Value = 28.828519197673
$ws.Cells.Item(17, 2).Value = 11.80527820367828
$ws.Cells.Item(17, 3).Value = 7.347161930494182
$ws.Cells.Item(17, 4).Value = 15.29164064605943
$ws.Cells.Item(17, 5).Value = 16.71858841299705
$ws.Cells.Item(17, 7).Value = 3.680610740331923
$ws.Cells.Item(17, 9).Value = 25.59455639492841
$ws.Cells.Item(17, 10).Value = 9.399690251570449
$ws.Cells.Item(17, 11).Value = 12.11134117754344
$ws.Cells.Item(17, 15).Value = 28.84887946354407
$ws.Cells.Item(18, 2).Value = 11.73919158180251
$ws.Cells.Item(18, 3).Value = 7.302965718818657
$ws.Cells.Item(18, 4).Value = 15.27076232477913
$ws.Cells.Item(18, 5).Value = 16.69731045244664
$ws.Cells.Item(18, 7).Value = 3.681118866014219
$ws.Cells.Item(18, 9).Value = 25.61266894611461
$ws.Cells.Item(18, 10).Value = 9.399917032486078
$ws.Cells.Item(18, 11).Value = 12.06593004958062
$ws.Cells.Item(18, 15).Value = 28.86118525833803
$ws.Cells.Item(19, 2).Value = 11.71676095591777
$ws.Cells.Item(19, 3).Value = 7.28796037126318
$ws.Cells.Item(19, 4).Value = 15.26374343013705
$ws.Cells.Item(19, 5).Value = 16.69016224057252
$ws.Cells.Item(19, 7).Value = 3.681292095089714
$ws.Cells.Item(19, 9).Value = 25.61888358069497
$ws.Cells.Item(19, 10).Value = 9.400017685837792
$ws.Cells.Item(19, 11).Value = 12.05054793237256
$ws.Cells.Item(19, 15).Value = 28.86545391094593
$ws.Cells.Item(20, 2).Value = 11.81748275027743
$ws.Cells.Item(20, 3).Value = 7.355321682605277
$ws.Cells.Item(20, 4).Value = 15.29552845897078
$ws.Cells.Item(20, 5).Value = 16.72255303867869
$ws.Cells.Item(20, 7).Value = 3.680517260844827
$ws.Cells.Item(20, 9).Value = 25.59124315133201
$ws.Cells.Item(20, 10).Value = 9.399659613330758
$ws.Cells.Item(20, 11).Value = 12.11974222865718
$ws.Cells.Item(20, 15).Value = 28.84665046607943
$ws.Cells.Item(21, 2).Value = 12.15174977249047
$ws.Cells.Item(21, 3).Value = 7.578556672536088
$ws.Cells.Item(21, 4).Value = 15.40573963946376
$ws.Cells.Item(21, 5).Value = 16.83521726540108
$ws.Cells.Item(21, 7).Value = 3.677996220313467
$ws.Cells.Item(21, 9).Value = 25.50409158191085
$ws.Cells.Item(21, 10).Value = 9.400114441215152
$ws.Cells.Item(21, 11).Value = 12.35154976179449
$ws.Cells.Item(21, 15).Value = 28.79059002949682
$ws.Cells.Item(22, 2).Value = 12.36675220777914
$ws.Cells.Item(22, 3).Value = 7.721905645845268
$ws.Cells.Item(22, 4).Value = 15.48023105026057
$ws.Cells.Item(22, 5).Value = 16.91162391516815
$ws.Cells.Item(22, 7).Value = 3.676409631648926
$ws.Cells.Item(22, 9).Value = 25.45140489077476
$ws.Cells.Item(22, 10).Value = 9.401644068249857
$ws.Cells.Item(22, 11).Value = 12.50230295301295
$ws.Cells.Item(22, 15).Value = 28.75926927292727
$ws.Cells.Item(23, 2).Value = 12.25234094040163
$ws.Cells.Item(23, 3).Value = 7.645645825169808
$ws.Cells.Item(23, 4).Value = 15.44025390387617
$ws.Cells.Item(23, 5).Value = 16.87059598861889
$ws.Cells.Item(23, 7).Value = 3.677250855043952
$ws.Cells.Item(23, 9).Value = 25.47913316719083
$ws.Cells.Item(23, 10).Value = 9.400715141752194
$ws.Cells.Item(23, 11).Value = 12.42192674850813
$ws.Cells.Item(23, 15).Value = 28.77549830380915
$ws.Cells.Item(24, 2).Value = 11.81196619850746
$ws.Cells.Item(24, 3).Value = 7.35163349344325
$ws.Cells.Item(24, 4).Value = 15.29376990672095
$ws.Cells.Item(24, 5).Value = 16.72075965119455
$ws.Cells.Item(24, 7).Value = 3.680559500722518
$ws.Cells.Item(24, 9).Value = 25.592739554959
$ws.Cells.Item(24, 10).Value = 9.399673032124889
$ws.Cells.Item(24, 11).Value = 12.11594432180985
$ws.Cells.Item(24, 15).Value = 28.84765632575028
$ws.Cells.Item(25, 2).Value = 11.32457102875365
$ws.Cells.Item(25, 3).Value = 7.025194103123983
$ws.Cells.Item(25, 4).Value = 15.14687586751298
$ws.Cells.Item(25, 5).Value = 16.57158375756452
$ws.Cells.Item(25, 7).Value = 3.684392378059873
$ws.Cells.Item(25, 9).Value = 25.73354979375063
$ws.Cells.Item(25, 10).Value = 9.403860005482887
$ws.Cells.Item(25, 11).Value = 11.78427963935967
$ws.Cells.Item(25, 15).Value = 28.94822905914208
